$wb = $excel.ActiveWorkbook

# The new "Spain" market sheet is modeled on the existing "Italy" sheet, so
# duplicate it (right after "Italy") to pick up all formatting, merged
# cells, column widths and styles, then rename the copy.
$italy = $wb.Worksheets.Item("Italy")
[void]$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# Fill in the market-specific values for Spain.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2061/T2049"

# Match the slightly different column widths used on the new sheet.
$spain.Columns.Item(1).ColumnWidth = 29.166666666666668
$spain.Columns.Item(2).ColumnWidth = 27.944010416666668
$spain.Columns.Item(3).ColumnWidth = 11.944010416666666
$spain.Columns.Item(4).ColumnWidth = 22.276041666666668

# "Italy" is no longer the active tab, so its remembered selection becomes
# the whole used range instead of the old single-cell selection.
[void]$italy.Range("A1:D11").Select()

# The newly added "Spain" sheet becomes the active tab/selection.
[void]$spain.Activate()
[void]$spain.Range("B10").Select()
